# This script applies the data corrections described in the commit message
# ("modified data (because there were some logic problems)") to the
# "données04" sheet. Columns A (a ratio/percentage-like value) and C (a
# count) are corrected on a handful of rows; column B is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 5.71
$ws.Range("C13").Value = 57

$ws.Range("A26").Value = 6.34
$ws.Range("C26").Value = 83

$ws.Range("A28").Value = 9.4499999999999993
$ws.Range("C28").Value = 92

$ws.Range("A30").Value = 5.09
$ws.Range("C30").Value = 95

$ws.Range("A31").Value = 15
$ws.Range("C31").Value = 70

$ws.Range("A32").Value = 21.85
$ws.Range("C32").Value = 94

$ws.Range("A33").Value = 17.89
$ws.Range("C33").Value = 71

$ws.Range("A35").Value = 9.48
$ws.Range("C35").Value = 65

$ws.Range("A36").Value = 7.37
$ws.Range("C36").Value = 95

$ws.Range("A38").Value = 60.089999999999996
$ws.Range("C38").Value = 81

$ws.Range("A40").Value = 9.0399999999999991
$ws.Range("C40").Value = 93

$ws.Range("A43").Value = 61.77
$ws.Range("C43").Value = 85

$ws.Range("A46").Value = 42.699999999999996
$ws.Range("C46").Value = 81

$ws.Range("A47").Value = 56.74
$ws.Range("C47").Value = 91

$ws.Range("A50").Value = 45.800000000000004
$ws.Range("C50").Value = 86

$ws.Range("A51").Value = 85.36
$ws.Range("C51").Value = 93

$ws.Range("A52").Value = 25.069999999999997
$ws.Range("C52").Value = 72
